$d = $word.ActiveDocument

$replacements = @(
    @("815÷4=", "533÷9="),
    @("459÷8=", "336÷9="),
    @("345÷6=", "870÷6="),
    @("494÷8=", "230÷3="),
    @("701÷6=", "248÷3="),
    @("748÷2=", "111÷7="),
    @("613÷3=", "169÷5="),
    @("264÷7=", "642÷4="),
    @("837÷9=", "275÷4="),
    @("915÷2=", "924÷4="),
    @("888÷6=", "915÷9="),
    @("446÷4=", "294÷8="),
    @("579÷3=", "408÷3="),
    @("695÷5=", "369÷8="),
    @("787÷7=", "767÷5="),
    @("414÷2=", "914÷8="),
    @("881÷7=", "667÷6="),
    @("348÷6=", "692÷5="),
    @("264÷3=", "282÷3="),
    @("961÷7=", "228÷7="),
    @("343÷7=", "856÷9="),
    @("658÷5=", "302÷7="),
    @("108÷2=", "441÷3="),
    @("251÷2=", "527÷4="),
    @("747÷7=", "300÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
